# Elimina participante desde panel de administración
# Removes the participant row "Alexis Sharon_20251130_195228" (row 10),
# shifting all subsequent participant rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()
